$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08230999999999999
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("M2").Value = 1.815761
$ws.Range("N2").Value = 5.447283000000001
$ws.Range("O2").Value = 0.07007596730428067
$ws.Range("P2").Value = 0.07007596730428067
$ws.Range("Q2").Value = 0.04981842930333334
$ws.Range("R2").Value = 0.44836586373
$ws.Range("S2").Value = 0.0005161995507752223
$ws.Range("T2").Value = 0.0005161995507752223
$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08230999999999999
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.5079540516959071
$ws.Range("P3").Value = 0.5079540516959072
$ws.Range("Q3").Value = 0.3611148584488889
$ws.Range("R3").Value = 3.25003372604
$ws.Range("S3").Value = 0.003741734340410085
$ws.Range("T3").Value = 0.003741734340410086
$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08230999999999999
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("M4").Value = 9.711409333333334
$ws.Range("N4").Value = 29.134228
$ws.Range("O4").Value = 0.3747940411327002
$ws.Range("P4").Value = 0.3747940411327002
$ws.Range("Q4").Value = 0.2664487007422222
$ws.Range("R4").Value = 2.39803830668
$ws.Range("S4").Value = 0.002760839744471308
$ws.Range("T4").Value = 0.002760839744471308
$ws.Range("G5").Value = 0.02743666666666666
$ws.Range("H5").Value = 0.08230999999999999
$ws.Range("I5").Value = 0.007366285056527356
$ws.Range("J5").Value = 0.007366285056527356
$ws.Range("M5").Value = 1.222391
$ws.Range("N5").Value = 3.667173
$ws.Range("O5").Value = 0.04717593986711188
$ws.Range("P5").Value = 0.04717593986711189
$ws.Range("Q5").Value = 0.03353833440333333
$ws.Range("R5").Value = 0.30184500963
$ws.Range("S5").Value = 0.0003475114208707394
$ws.Range("T5").Value = 0.0003475114208707395
$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("M6").Value = 1.815761
$ws.Range("N6").Value = 5.447283000000001
$ws.Range("O6").Value = 0.07007596730428067
$ws.Range("P6").Value = 0.07007596730428067
$ws.Range("Q6").Value = 6.116080433369
$ws.Range("R6").Value = 55.044723900321
$ws.Range("S6").Value = 0.06337249119170771
$ws.Range("T6").Value = 0.06337249119170771
$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("O7").Value = 0.5079540516959071
$ws.Range("P7").Value = 0.5079540516959072
$ws.Range("Q7").Value = 44.33314239014533
$ws.Range("R7").Value = 398.998281511308
$ws.Range("S7").Value = 0.4593631012914285
$ws.Range("T7").Value = 0.4593631012914286
$ws.Range("G8").Value = 3.368329
$ws.Range("I8").Value = 0.9043398704228307
$ws.Range("J8").Value = 0.9043398704228307
$ws.Range("M8").Value = 9.711409333333334
$ws.Range("N8").Value = 29.134228
$ws.Range("O8").Value = 0.3747940411327002
$ws.Range("P8").Value = 0.3747940411327002
$ws.Range("Q8").Value = 32.71122168833733
$ws.Range("R8").Value = 294.400995195036
$ws.Range("S8").Value = 0.3389411945931952
$ws.Range("T8").Value = 0.3389411945931952
$ws.Range("G9").Value = 3.368329
$ws.Range("I9").Value = 0.9043398704228307
$ws.Range("J9").Value = 0.9043398704228307
$ws.Range("M9").Value = 1.222391
$ws.Range("N9").Value = 3.667173
$ws.Range("O9").Value = 0.04717593986711188
$ws.Range("P9").Value = 0.04717593986711189
$ws.Range("Q9").Value = 4.117415054638999
$ws.Range("R9").Value = 37.056735491751
$ws.Range("S9").Value = 0.04266308334649922
$ws.Range("T9").Value = 0.04266308334649922
$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("M10").Value = 1.815761
$ws.Range("N10").Value = 5.447283000000001
$ws.Range("O10").Value = 0.07007596730428067
$ws.Range("P10").Value = 0.07007596730428067
$ws.Range("Q10").Value = 0.5971341887283333
$ws.Range("R10").Value = 5.374207698555
$ws.Range("S10").Value = 0.006187276561797749
$ws.Range("T10").Value = 0.006187276561797749
$ws.Range("G11").Value = 0.3288616666666667
$ws.Range("H11").Value = 0.9865849999999999
$ws.Range("I11").Value = 0.08829384452064198
$ws.Range("J11").Value = 0.08829384452064198
$ws.Range("O11").Value = 0.5079540516959071
$ws.Range("P11").Value = 0.5079540516959072
$ws.Range("Q11").Value = 4.328398768348888
$ws.Range("R11").Value = 38.95558891514
$ws.Range("S11").Value = 0.04484921606406857
$ws.Range("T11").Value = 0.04484921606406857
$ws.Range("G12").Value = 0.3288616666666667
$ws.Range("H12").Value = 0.9865849999999999
$ws.Range("I12").Value = 0.08829384452064198
$ws.Range("J12").Value = 0.08829384452064198
$ws.Range("M12").Value = 9.711409333333334
$ws.Range("N12").Value = 29.134228
$ws.Range("O12").Value = 0.3747940411327002
$ws.Range("P12").Value = 0.3747940411327002
$ws.Range("Q12").Value = 3.193710259042223
$ws.Range("R12").Value = 28.74339233138
$ws.Range("S12").Value = 0.03309200679503373
$ws.Range("T12").Value = 0.03309200679503373
$ws.Range("G13").Value = 0.3288616666666667
$ws.Range("H13").Value = 0.9865849999999999
$ws.Range("I13").Value = 0.08829384452064198
$ws.Range("J13").Value = 0.08829384452064198
$ws.Range("M13").Value = 1.222391
$ws.Range("N13").Value = 3.667173
$ws.Range("O13").Value = 0.04717593986711188
$ws.Range("P13").Value = 0.04717593986711189
$ws.Range("Q13").Value = 0.4019975415783333
$ws.Range("R13").Value = 3.617977874205
$ws.Range("S13").Value = 0.004165345099741933
$ws.Range("T13").Value = 0.004165345099741933
